$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 321, pushing the existing rows 321:398 down to 324:401
$ws.Rows("321:323").Insert()

# Populate the 3 newly inserted rows with the new weekly price entries
# (Nectarín, Super Queen variety, date 2023-01-06 = serial 44932)

# Row 321: Super Queen / Especial
$ws.Range("A321").Value = 11
$ws.Range("B321").Value = "Vega Monumental Concepción"
$ws.Range("C321").Value = "Bíobío"
$ws.Range("D321").Value = 44932
$ws.Range("D321").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E321").Value = 8
$ws.Range("F321").Value = "Fruta"
$ws.Range("G321").Value = 100103
$ws.Range("H321").Value = "Frutos de hueso (carozo)"
$ws.Range("I321").Value = 100103006
$ws.Range("J321").Value = "Nectarín"
$ws.Range("K321").Value = "Super Queen"
$ws.Range("L321").Value = "Especial"
$ws.Range("M321").Value = 50
$ws.Range("N321").Value = 14000
$ws.Range("O321").Value = 14000
$ws.Range("P321").Value = 14000
$ws.Range("Q321").Value = "$/caja 15 kilos empedrada"
$ws.Range("R321").Value = "Región de O'Higgins"
$ws.Range("S321").Value = 933
$ws.Range("T321").Value = 15

# Row 322: Super Queen / Primera
$ws.Range("A322").Value = 11
$ws.Range("B322").Value = "Vega Monumental Concepción"
$ws.Range("C322").Value = "Bíobío"
$ws.Range("D322").Value = 44932
$ws.Range("D322").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E322").Value = 8
$ws.Range("F322").Value = "Fruta"
$ws.Range("G322").Value = 100103
$ws.Range("H322").Value = "Frutos de hueso (carozo)"
$ws.Range("I322").Value = 100103006
$ws.Range("J322").Value = "Nectarín"
$ws.Range("K322").Value = "Super Queen"
$ws.Range("L322").Value = "Primera"
$ws.Range("M322").Value = 50
$ws.Range("N322").Value = 12000
$ws.Range("O322").Value = 12000
$ws.Range("P322").Value = 12000
$ws.Range("Q322").Value = "$/caja 15 kilos empedrada"
$ws.Range("R322").Value = "Región de O'Higgins"
$ws.Range("S322").Value = 800
$ws.Range("T322").Value = 15

# Row 323: Super Queen / Segunda
$ws.Range("A323").Value = 11
$ws.Range("B323").Value = "Vega Monumental Concepción"
$ws.Range("C323").Value = "Bíobío"
$ws.Range("D323").Value = 44932
$ws.Range("D323").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E323").Value = 8
$ws.Range("F323").Value = "Fruta"
$ws.Range("G323").Value = 100103
$ws.Range("H323").Value = "Frutos de hueso (carozo)"
$ws.Range("I323").Value = 100103006
$ws.Range("J323").Value = "Nectarín"
$ws.Range("K323").Value = "Super Queen"
$ws.Range("L323").Value = "Segunda"
$ws.Range("M323").Value = 50
$ws.Range("N323").Value = 10000
$ws.Range("O323").Value = 10000
$ws.Range("P323").Value = 10000
$ws.Range("Q323").Value = "$/caja 15 kilos empedrada"
$ws.Range("R323").Value = "Región de O'Higgins"
$ws.Range("S323").Value = 667
$ws.Range("T323").Value = 15
